$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Administration" paragraph: colour the whole paragraph (run + para mark)
#    red (FF0000).
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Administration", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Paragraphs(1).Range.Font.Color = 255
}

# ---------------------------------------------------------------------------
# 2. "Mail form " run + following space run get coloured red (but not the
#    line break or the italic text that follows).
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Mail form ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Font.Color = 255
}

# ---------------------------------------------------------------------------
# 3. "Filter by price, country , ... " - the run split caused by the
#    now-removed _GoBack bookmark is re-merged into a single run. The
#    visible text does not change, only the run/bookmark structure, so we
#    round-trip the text through a placeholder to force Word to rebuild a
#    single run without the bookmark.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute(" by price, country , … ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Text = "`u{E000}PLACEHOLDER`u{E000}"
    $r2 = $d.Content
    $r2.Find.Execute("`u{E000}PLACEHOLDER`u{E000}", $true, $false, $false, $false, $false, $true, 1, $false, " by price, country , … ", 2)
}

# ---------------------------------------------------------------------------
# 4. "Map " run + following space run get coloured red.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Map ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Font.Color = 255
}

# ---------------------------------------------------------------------------
# 5. "Calendar" run loses its underline (the <w:u> element is removed
#    entirely, it is not merely toggled to "none"). The high level
#    Font.Underline API always serialises an explicit value, so we rebuild
#    just this run via InsertXML after clearing its text, which lets us
#    omit the underline property altogether.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Calendar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $startPos = $r.Start
    $r.Text = ""
    $collapsed = $d.Range($startPos, $startPos)
    $runXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic" w:cs="Angsana New"/><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>Calendar</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $collapsed.InsertXML($runXml)
}

# ---------------------------------------------------------------------------
# 6. "Price convert " run + following space run get coloured red.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Price convert", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Font.Color = 255
}

# ---------------------------------------------------------------------------
# 7. "Use bank API and convert price..." run is split in two right after
#    "Use bank API a", with a (collapsed) _GoBack bookmark marking the split
#    point, mirroring where the author's cursor was when the file was last
#    saved.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Use bank API a", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPos = $r.End
    $splitPoint = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_GoBack", $splitPoint)
}
